$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.1636203333333333
$ws.Range("H2").Value = 0.490861
$ws.Range("I2").Value = 0.3345941539187231
$ws.Range("J2").Value = 0.3345941539187231
$ws.Range("M2").Value = 0.4067693333333334
$ws.Range("N2").Value = 1.220308
$ws.Range("O2").Value = 0.1164607724076721
$ws.Range("P2").Value = 0.1164607724076721
$ws.Range("Q2").Value = 0.06655573390977779
$ws.Range("R2").Value = 0.5990016051880001
$ws.Range("S2").Value = 0.03896709360846601
$ws.Range("T2").Value = 0.03896709360846601
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.1636203333333333
$ws.Range("H3").Value = 0.490861
$ws.Range("I3").Value = 0.3345941539187231
$ws.Range("J3").Value = 0.3345941539187231
$ws.Range("O3").Value = 0.3361456919197101
$ws.Range("P3").Value = 0.33614569191971
$ws.Range("Q3").Value = 0.1921026519385556
$ws.Range("R3").Value = 1.728923867447
$ws.Range("S3").Value = 0.1124723833812992
$ws.Range("T3").Value = 0.1124723833812991
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.1636203333333333
$ws.Range("H4").Value = 0.490861
$ws.Range("I4").Value = 0.3345941539187231
$ws.Range("J4").Value = 0.3345941539187231
$ws.Range("M4").Value = 1.368657333333333
$ws.Range("N4").Value = 4.105972
$ws.Range("O4").Value = 0.391855720526518
$ws.Range("P4").Value = 0.391855720526518
$ws.Range("Q4").Value = 0.2239401690991111
$ws.Range("R4").Value = 2.015461521892
$ws.Range("S4").Value = 0.1311126332677819
$ws.Range("T4").Value = 0.1311126332677819
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.1636203333333333
$ws.Range("H5").Value = 0.490861
$ws.Range("I5").Value = 0.3345941539187231
$ws.Range("J5").Value = 0.3345941539187231
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.543256
$ws.Range("N5").Value = 1.629768
$ws.Range("O5").Value = 0.1555378151460999
$ws.Range("P5").Value = 0.1555378151460999
$ws.Range("Q5").Value = 0.08888772780533333
$ws.Range("R5").Value = 0.7999895502479999
$ws.Range("S5").Value = 0.05204204366117605
$ws.Range("T5").Value = 0.05204204366117605
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.325391
$ws.Range("H6").Value = 0.9761730000000001
$ws.Range("I6").Value = 0.665405846081277
$ws.Range("J6").Value = 0.665405846081277
$ws.Range("M6").Value = 0.4067693333333334
$ws.Range("N6").Value = 1.220308
$ws.Range("O6").Value = 0.1164607724076721
$ws.Range("P6").Value = 0.1164607724076721
$ws.Range("Q6").Value = 0.1323590801426667
$ws.Range("R6").Value = 1.191231721284
$ws.Range("S6").Value = 0.07749367879920607
$ws.Range("T6").Value = 0.07749367879920606
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.325391
$ws.Range("H7").Value = 0.9761730000000001
$ws.Range("I7").Value = 0.665405846081277
$ws.Range("J7").Value = 0.665405846081277
$ws.Range("O7").Value = 0.3361456919197101
$ws.Range("P7").Value = 0.33614569191971
$ws.Range("Q7").Value = 0.3820336552523334
$ws.Range("R7").Value = 3.438302897271
$ws.Range("S7").Value = 0.223673308538411
$ws.Range("T7").Value = 0.2236733085384109
$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.325391
$ws.Range("H8").Value = 0.9761730000000001
$ws.Range("I8").Value = 0.665405846081277
$ws.Range("J8").Value = 0.665405846081277
$ws.Range("M8").Value = 1.368657333333333
$ws.Range("N8").Value = 4.105972
$ws.Range("O8").Value = 0.391855720526518
$ws.Range("P8").Value = 0.391855720526518
$ws.Range("Q8").Value = 0.4453487783506667
$ws.Range("R8").Value = 4.008139005156001
$ws.Range("S8").Value = 0.2607430872587361
$ws.Range("T8").Value = 0.2607430872587361
$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.325391
$ws.Range("H9").Value = 0.9761730000000001
$ws.Range("I9").Value = 0.665405846081277
$ws.Range("J9").Value = 0.665405846081277
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.543256
$ws.Range("N9").Value = 1.629768
$ws.Range("O9").Value = 0.1555378151460999
$ws.Range("P9").Value = 0.1555378151460999
$ws.Range("Q9").Value = 0.176770613096
$ws.Range("R9").Value = 1.590935517864
$ws.Range("S9").Value = 0.1034957714849239
$ws.Range("T9").Value = 0.1034957714849238
